$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 3 <-> Row 6 swap species/observation data (and rows 4,5,7 rotate the
# record "Id" + coordinates amongst themselves), matching the source diff.
# ---------------------------------------------------------------------------

# --- Row 3 (was the "Tretaig hackspett" record, becomes the "Garnlav" one) ---
$ws.Range("A3").Value = 111741082
$ws.Range("B3").Value = 77515
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = "Garnlav"
$ws.Range("G3").Value = "Alectoria sarmentosa"
$ws.Range("H3").Value = "(Ach.) Ach."
$ws.Range("L3").Value = $null
$ws.Range("M3").Value = $null
$ws.Range("Q3").Value = 331469
$ws.Range("R3").Value = 6627064
$ws.Range("Z3").Value = $null
$ws.Range("AB3").Value = $null

# --- Row 4 (Id rotates; same species data) ---
$ws.Range("A4").Value = 111741025
$ws.Range("Q4").Value = 331437
$ws.Range("R4").Value = 6627065
$ws.Range("Z4").Value = $null
$ws.Range("AB4").Value = $null

# --- Row 5 (Id rotates; same species data) ---
$ws.Range("A5").Value = 111741014
$ws.Range("Q5").Value = 331429
$ws.Range("R5").Value = 6627058
$ws.Range("Z5").Value = $null
$ws.Range("AB5").Value = $null

# --- Row 6 (was the "Garnlav" record, becomes the "Tretaig hackspett" one) ---
$ws.Range("A6").Value = 111741120
$ws.Range("B6").Value = 56398
$ws.Range("E6").Value = 100109
$ws.Range("F6").Value = "Tretåig hackspett"
$ws.Range("G6").Value = "Picoides tridactylus"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("J6").Value = $null
$ws.Range("M6").Value = "färska spår"
$ws.Range("Q6").Value = 331469
$ws.Range("R6").Value = 6627064
$ws.Range("Z6").Value = $null
$ws.Range("AB6").Value = $null
$ws.Range("AF6").Value = $null

# --- Row 7 (Id rotates; same species data) ---
$ws.Range("A7").Value = 111741038
$ws.Range("Q7").Value = 331443
$ws.Range("R7").Value = 6627065
$ws.Range("Z7").Value = $null
$ws.Range("AB7").Value = $null

# --- Row 8 (coordinates rounded; start/end time columns dropped) ---
$ws.Range("Q8").Value = 331462
$ws.Range("R8").Value = 6627063
$ws.Range("Z8").Value = $null
$ws.Range("AB8").Value = $null
